{"js": "// Applies the resume content updates described in the commit:\n//  1. \"21 years\" -> \"15+ years\" in the professional summary\n//  2. FLEEM bullet gets more detail about Twilio/predictive dialer usage\n//  3. Salsa Labs bullets replaced with geospatial/CRM-focused bullets\n//  4. Praxis Project bullets replaced with expanded leadership/training bullets\n//  5. New trailing bullet added to Lake Research Partners section\n//  6. New trailing bullet added to Feldman Group section\n\nconst body = context.document.body;\n\n// Helper: find the (single) paragraph whose text matches `needle` exactly\n// by locating it via a unique substring search, then returning its\n// paragraph object.\nasync function findParagraph(needle) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Paragraph not found for search text: \" + needle);\n  }\n  return results.items[0].paragraphs.getFirst();\n}\n\n// Helper: replace a run of consecutive bullet paragraphs (identified by\n// their current text) with a new list of bullet strings. The first\n// paragraph's text is overwritten in place (preserving its paragraph\n// identity/formatting); the remaining new bullets are inserted after it;\n// then any leftover old paragraphs (beyond the number of new bullets) are\n// deleted.\nasync function replaceBulletRun(oldTexts, newTexts) {\n  // Resolve every old paragraph up front (search the original text before\n  // any edits in this run are made).\n  const oldParagraphs = [];\n  for (const t of oldTexts) {\n    oldParagraphs.push(await findParagraph(t));\n  }\n\n  // Overwrite the first paragraph with the first new bullet.\n  oldParagraphs[0].insertText(newTexts[0], \"Replace\");\n  await context.sync();\n\n  // Insert the remaining new bullets, in order, right after the first.\n  let anchor = oldParagraphs[0];\n  for (let i = 1; i < newTexts.length; i++) {\n    anchor = anchor.insertParagraph(newTexts[i], \"After\");\n    await context.sync();\n  }\n\n  // Delete whatever old paragraphs are left over (everything beyond the\n  // first, which already got reused above).\n  for (let i = 1; i < oldParagraphs.length; i++) {\n    oldParagraphs[i].delete();\n  }\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 1. Update years of experience in the professional summary\n// ---------------------------------------------------------------------\nconst summaryMatches = body.search(\"21 years\", { matchCase: true });\nsummaryMatches.load(\"text\");\nawait context.sync();\nsummaryMatches.items[0].insertText(\"15+ years\", \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2. Expand the FLEEM / Twilio bullet (Progressive Change Campaign Committee)\n// ---------------------------------------------------------------------\nconst fleemPar = await findParagraph(\n  \"Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\"\n);\nfleemPar.insertText(\n  \"\\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys\",\n  \"Replace\"\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3. Replace the Salsa Labs, Inc. bullets\n// ---------------------------------------------------------------------\nawait replaceBulletRun(\n  [\n    \"Developed software solutions for political campaigns and advocacy groups\",\n    \"Built web applications for voter engagement and campaign management\",\n    \"Integrated third-party APIs and data sources for campaign tools\",\n    \"Collaborated with political strategists to translate requirements into technical solutions\",\n  ],\n  [\n    \"\\u2022 Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously\",\n    \"\\u2022 Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers\",\n    \"\\u2022 Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\",\n    \"\\u2022 Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs\",\n    \"\\u2022 Collaborated with political strategists to translate geospatial requirements into technical solutions\",\n    \"\\u2022 Handled billions of records with millions of columns in high-performance CRM system\",\n  ]\n);\n\n// ---------------------------------------------------------------------\n// 4. Replace The Praxis Project bullets\n// ---------------------------------------------------------------------\nawait replaceBulletRun(\n  [\n    \"Integrated technology solutions within organizational frameworks for social justice organizations\",\n    \"Developed data management systems for community organizing efforts\",\n    \"Provided technical training and support to nonprofit staff\",\n    \"Built custom applications for community engagement and advocacy\",\n  ],\n  [\n    \"\\u2022 Led technology operations for multi-million dollar organization while assisting in search for full-time CTO\",\n    \"\\u2022 Directed all technology decisions and practices for massive multinational non-governmental organization\",\n    \"\\u2022 Developed comprehensive frameworks for internal and external technology audits\",\n    \"\\u2022 Led training initiatives for beneficiaries on spatial and Census data analysis for public health research\",\n    \"\\u2022 Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL\",\n    \"\\u2022 Managed technology infrastructure supporting community health initiatives across multiple countries\",\n    \"\\u2022 Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation\",\n  ]\n);\n\n// ---------------------------------------------------------------------\n// 5. Add a new trailing bullet to the Lake Research Partners section\n// ---------------------------------------------------------------------\nconst lakeResearchLastBullet = await findParagraph(\n  \"Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\"\n);\nlakeResearchLastBullet.insertParagraph(\n  \"\\u2022 Trained staff on building Python tooling for report generation and analysis\",\n  \"After\"\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 6. Add a new trailing bullet to The Feldman Group section\n// ---------------------------------------------------------------------\nconst feldmanLastBullet = await findParagraph(\n  \"Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\"\n);\nfeldmanLastBullet.insertParagraph(\n  \"\\u2022 Trained staff on PHP/MySQL for data analysis and reporting systems\",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "# Applies the resume content updates described in the commit:\n#  1. \"21 years\" -> \"15+ years\" in the professional summary\n#  2. FLEEM bullet gets more detail about Twilio/predictive dialer usage\n#  3. Salsa Labs bullets replaced with geospatial/CRM-focused bullets\n#  4. Praxis Project bullets replaced with expanded leadership/training bullets\n#  5. New trailing bullet added to Lake Research Partners section\n#  6. New trailing bullet added to Feldman Group section\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($doc, $needle) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -like \"*$needle*\") {\n            return $p\n        }\n    }\n    return $null\n}\n\n# Replaces a contiguous run of bullet paragraphs (identified by their\n# current, unique text) with a new set of bullet strings. The first\n# paragraph in the run is reused (its text overwritten) and the remaining\n# new bullets are inserted right after it; any leftover old paragraphs are\n# deleted FIRST (from last to first) so the live Paragraph references for\n# the surviving paragraphs do not get invalidated by the later inserts.\nfunction Replace-BulletRun($doc, $oldTexts, $newTexts) {\n    $oldParagraphs = @()\n    foreach ($t in $oldTexts) {\n        $oldParagraphs += , (Find-ParagraphByText $doc $t)\n    }\n\n    # Delete the leftover old paragraphs (everything after the first) from\n    # last to first before making any other edits.\n    for ($i = $oldParagraphs.Count - 1; $i -ge 1; $i--) {\n        $oldParagraphs[$i].Range.Delete()\n    }\n\n    # Overwrite the first paragraph with the first new bullet.\n    $first = $oldParagraphs[0]\n    $first.Range.Text = $newTexts[0]\n\n    # Insert the remaining new bullets, in order, right after the first.\n    $anchor = $first\n    for ($i = 1; $i -lt $newTexts.Count; $i++) {\n        $anchor.Range.InsertParagraphAfter()\n        $anchor = $anchor.Next()\n        $anchor.Range.Text = $newTexts[$i]\n    }\n}\n\n# ---------------------------------------------------------------------\n# 1. Update years of experience in the professional summary\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"21 years\"\n$find.Replacement.Text = \"15+ years\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# ---------------------------------------------------------------------\n# 2. Expand the FLEEM / Twilio bullet (Progressive Change Campaign Committee)\n# ---------------------------------------------------------------------\n$fleemPar = Find-ParagraphByText $d \"Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\"\n$fleemPar.Range.Text = \"\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys\"\n\n# ---------------------------------------------------------------------\n# 3. Replace the Salsa Labs, Inc. bullets\n# ---------------------------------------------------------------------\n$salsaOld = @(\n    \"Developed software solutions for political campaigns and advocacy groups\",\n    \"Built web applications for voter engagement and campaign management\",\n    \"Integrated third-party APIs and data sources for campaign tools\",\n    \"Collaborated with political strategists to translate requirements into technical solutions\"\n)\n$salsaNew = @(\n    \"\u2022 Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously\",\n    \"\u2022 Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers\",\n    \"\u2022 Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\",\n    \"\u2022 Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs\",\n    \"\u2022 Collaborated with political strategists to translate geospatial requirements into technical solutions\",\n    \"\u2022 Handled billions of records with millions of columns in high-performance CRM system\"\n)\nReplace-BulletRun $d $salsaOld $salsaNew\n\n# ---------------------------------------------------------------------\n# 4. Replace The Praxis Project bullets\n# ---------------------------------------------------------------------\n$praxisOld = @(\n    \"Integrated technology solutions within organizational frameworks for social justice organizations\",\n    \"Developed data management systems for community organizing efforts\",\n    \"Provided technical training and support to nonprofit staff\",\n    \"Built custom applications for community engagement and advocacy\"\n)\n$praxisNew = @(\n    \"\u2022 Led technology operations for multi-million dollar organization while assisting in search for full-time CTO\",\n    \"\u2022 Directed all technology decisions and practices for massive multinational non-governmental organization\",\n    \"\u2022 Developed comprehensive frameworks for internal and external technology audits\",\n    \"\u2022 Led training initiatives for beneficiaries on spatial and Census data analysis for public health research\",\n    \"\u2022 Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL\",\n    \"\u2022 Managed technology infrastructure supporting community health initiatives across multiple countries\",\n    \"\u2022 Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation\"\n)\nReplace-BulletRun $d $praxisOld $praxisNew\n\n# ---------------------------------------------------------------------\n# 5. Add a new trailing bullet to the Lake Research Partners section\n# ---------------------------------------------------------------------\n$lakeResearchLastBullet = Find-ParagraphByText $d \"Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\"\n$lakeResearchLastBullet.Range.InsertParagraphAfter()\n$newPar = $lakeResearchLastBullet.Next()\n$newPar.Range.Text = \"\u2022 Trained staff on building Python tooling for report generation and analysis\"\n\n# ---------------------------------------------------------------------\n# 6. Add a new trailing bullet to The Feldman Group section\n# ---------------------------------------------------------------------\n$feldmanLastBullet = Find-ParagraphByText $d \"Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\"\n$feldmanLastBullet.Range.InsertParagraphAfter()\n$newPar2 = $feldmanLastBullet.Next()\n$newPar2.Range.Text = \"\u2022 Trained staff on PHP/MySQL for data analysis and reporting systems\"\n"}
